$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update B15/C15 and B17/C17 to "Si", matching the style of B13/C13 (s=11)
$ws.Range("B15").Value = "Si"
$ws.Range("C15").Value = "Si"
$ws.Range("B15:C15").Style = $ws.Range("B13:C13").Style

$ws.Range("B17").Value = "Si"
$ws.Range("C17").Value = "Si"
$ws.Range("B17:C17").Style = $ws.Range("B13:C13").Style

# Update the active selection to C11 (was F12)
$ws.Range("C11").Select()
